# Daily scrape update - 2025-07-28 03:49:00 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (stored XML width = ColumnWidth + 5/6) ---
$ws.Columns.Item(3).ColumnWidth = 45.166666666666664   # C: 37 -> 46
$ws.Columns.Item(4).ColumnWidth = 24.166666666666668   # D: 22 -> 25
$ws.Columns.Item(6).ColumnWidth = 16.166666666666668   # F: 15 -> 17

# --- Row 2: opportunity 1326581 -> 1326587 ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1326587"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1326587"
$ws.Range("C2").Value = "Digital Marketing"
$ws.Range("D2").Value = "2750 Cascais, Portugal"
$ws.Range("F2").Value = "9 applicants"
$ws.Range("H2").Value = "MQ Capital"

# --- Row 3: opportunity 1325142 -> 1326580 ---
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "1326580"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1326580"
$ws.Range("C3").Value = "International Sales Representative - Intern"
$ws.Range("F3").Value = "0 applicants"
$ws.Range("G3").Value = "6 - 18 Months"
$ws.Range("H3").Value = "Brand Corridor (Pvt) Ltd"

# --- Row 4 (new): opportunity 1323793 ---
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "1323793"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1323793"
$ws.Range("C4").Value = "Digital Marketing"
$ws.Range("D4").Value = "Faro, Portugal"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "137 applicants"
$ws.Range("G4").Value = "3 - 6 Months"
$ws.Range("H4").Value = "Bed and Sun"
